# Appends two new observation rows (9 and 10) to the Artfynd sheet,
# matching the "Garnlav" (Alectoria sarmentosa) records added upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a date/time-looking string as *text* (no auto date-serial
# conversion, no leftover cell style) by staging it in a scratch cell that is
# explicitly formatted as Text, copying it, and pasting *values only* into
# the real destination. The scratch cell is cleared immediately afterwards so
# it leaves no trace in the sheet's used range.
function Set-TextValue($rangeAddress, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($rangeAddress).PasteSpecial(-4163)   # -4163 = xlPasteValues
    $scratch.Clear()
}

# ---------------------------------------------------------------------
# Row 9
# ---------------------------------------------------------------------
$ws.Range("A9").Value = 131157730
$ws.Range("B9").Value = 79245
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 6425
$ws.Range("F9").Value = "Garnlav"
$ws.Range("G9").Value = "Alectoria sarmentosa"
$ws.Range("H9").Value = "(Ach.) Ach."
$ws.Range("P9").Value = "Valmyran, Valmyran, Ång"
$ws.Range("Q9").Value = 715048
$ws.Range("R9").Value = 7090771
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = "Västerbotten"
$ws.Range("U9").Value = "Bjurholm"
$ws.Range("V9").Value = "Ångermanland"
$ws.Range("W9").Value = "Bjurholm"
Set-TextValue "Y9" "2026-02-15"
$ws.Range("Z9").Value = "11:40"
Set-TextValue "AA9" "2026-02-15"
$ws.Range("AB9").Value = "11:40"
$ws.Range("AC9").Value = "På sälg"
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AW9").Value = "Anne Siivola"
$ws.Range("AX9").Value = "Anne Siivola"

# ---------------------------------------------------------------------
# Row 10
# ---------------------------------------------------------------------
$ws.Range("A10").Value = 131159001
$ws.Range("B10").Value = 79245
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("P10").Value = "Valmyran, Valmyran, Ång"
$ws.Range("Q10").Value = 714973
$ws.Range("R10").Value = 7090850
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = "Västerbotten"
$ws.Range("U10").Value = "Bjurholm"
$ws.Range("V10").Value = "Ångermanland"
$ws.Range("W10").Value = "Bjurholm"
Set-TextValue "Y10" "2026-02-15"
$ws.Range("Z10").Value = "13:30"
Set-TextValue "AA10" "2026-02-15"
$ws.Range("AB10").Value = "13:30"
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AW10").Value = "Anne Siivola"
$ws.Range("AX10").Value = "Anne Siivola"
